# "Generate Report for Handoff"
#
# Refreshes the handoff-report timestamps:
#  - zh-cn sheet: Test`1.md's "Latest Handoff Datetime" (H3) moves forward
#    to the latest handoff run.
#  - de-de sheet: Test`1.md and Test`2.md share the same "Latest Handoff
#    Datetime" value (H3 and H4); both advance together to the new run time.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2017-11-04 16:03:20"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2017-11-04 16:03:26"
$wsDeDe.Range("H4").Value = "2017-11-04 16:03:26"
